$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextCell "B2" 'Bitcoin'
Set-TextCell "C2" 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextCell "D2" '88.000.16'
Set-TextCell "E2" '  -2.51%  '

Set-TextCell "B3" 'Ethereum'
Set-TextCell "C3" 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextCell "D3" '3.068.21'
Set-TextCell "E3" '  -5.18%  '

Set-TextCell "B4" 'TetherUSD'
Set-TextCell "C4" 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextCell "D4" '1.00'
Set-TextCell "E4" '  -0.03%  '

Set-TextCell "B5" 'Solana'
Set-TextCell "C5" 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell "D5" '209.79'
Set-TextCell "E5" '  -4.70%  '

Set-TextCell "B6" 'BNB'
Set-TextCell "C6" 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextCell "D6" '619.45'
Set-TextCell "E6" '  -2.34%  '

Set-TextCell "B7" 'Dogecoin'
Set-TextCell "C7" 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell "D7" '0.371'
Set-TextCell "E7" '  -7.16%  '

Set-TextCell "B8" 'XRP'
Set-TextCell "C8" 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextCell "D8" '0.813'
Set-TextCell "E8" '  +15.56%  '

Set-TextCell "B9" 'USDC'
Set-TextCell "C9" 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextCell "D9" '1.00'
Set-TextCell "E9" '  +0.00%  '

Set-TextCell "B10" 'LidoStakedEther'
Set-TextCell "C10" 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
Set-TextCell "D10" '3.064.18'
Set-TextCell "E10" '  -5.20%  '

Set-TextCell "B11" 'Cardano'
Set-TextCell "C11" 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell "D11" '0.620'
Set-TextCell "E11" '  +6.02%  '

Set-TextCell "B12" 'TRON'
Set-TextCell "C12" 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell "D12" '0.179'
Set-TextCell "E12" '  -1.74%  '

Set-TextCell "B13" 'ShibaInu'
Set-TextCell "C13" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell "D13" '0.0000238'
Set-TextCell "E13" '  -11.07%  '

Set-TextCell "B14" 'Toncoin'
Set-TextCell "C14" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell "D14" '5.28'
Set-TextCell "E14" '  -3.58%  '

Set-TextCell "B15" 'WrappedBTC'
Set-TextCell "C15" 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell "D15" '87.808.82'
Set-TextCell "E15" '  -2.44%  '

Set-TextCell "B16" 'WrappedliquidstakedEther2.0'
Set-TextCell "C16" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell "D16" '3.631.04'
Set-TextCell "E16" '  -5.21%  '

Set-TextCell "B17" 'Avalanche'
Set-TextCell "C17" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell "D17" '31.84'
Set-TextCell "E17" '  -7.50%  '

Set-TextCell "B18" 'WrappedEther'
Set-TextCell "C18" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell "D18" '3.053.99'
Set-TextCell "E18" '  -5.37%  '

Set-TextCell "B19" 'SuiNetwork'
Set-TextCell "C19" 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextCell "D19" '3.21'
Set-TextCell "E19" '  -7.88%  '

Set-TextCell "B20" 'PEPE'
Set-TextCell "C20" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell "D20" '0.0000201'
Set-TextCell "E20" '  -13.68%  '

Set-TextCell "B21" 'Chainlink'
Set-TextCell "C21" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell "D21" '13.21'
Set-TextCell "E21" '  -4.28%  '

Set-TextCell "B22" 'BitcoinCash'
Set-TextCell "C22" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell "D22" '419.96'
Set-TextCell "E22" '  -5.69%  '

Set-TextCell "B23" 'Uniswap'
Set-TextCell "C23" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell "D23" '8.13'
Set-TextCell "E23" '  -7.56%  '

Set-TextCell "B24" 'Polkadot'
Set-TextCell "C24" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell "D24" '4.88'
Set-TextCell "E24" '  -6.00%  '

Set-TextCell "B25" 'NEARProtocol'
Set-TextCell "C25" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell "D25" '5.47'
Set-TextCell "E25" '  +2.87%  '

Set-TextCell "B26" 'Aptos'
Set-TextCell "C26" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell "D26" '11.75'
Set-TextCell "E26" '  -3.78%  '

Set-TextCell "B27" 'Litecoin'
Set-TextCell "C27" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell "D27" '81.86'
Set-TextCell "E27" '  -2.18%  '

Set-TextCell "B28" 'Dai'
Set-TextCell "C28" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell "D28" '0.999'
Set-TextCell "E28" '  -0.06%  '

Set-TextCell "B29" 'Binance-PegBSC-USD'
Set-TextCell "C29" 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell "D29" '1.07'
Set-TextCell "E29" '  +7.10%  '

Set-TextCell "B30" 'Cronos'
Set-TextCell "C30" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell "D30" '0.164'
Set-TextCell "E30" '  +0.40%  '

Set-TextCell "B31" 'InternetComputer(DFINITY)'
Set-TextCell "C31" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell "D31" '8.04'
Set-TextCell "E31" '  -7.44%  '

Set-TextCell "B32" 'Bittensor'
Set-TextCell "C32" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell "D32" '505.10'
Set-TextCell "E32" '  -9.05%  '

Set-TextCell "B33" 'dogwifhat'
Set-TextCell "C33" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell "D33" '3.54'
Set-TextCell "E33" '  -16.12%  '

Set-TextCell "B34" 'RenderToken'
Set-TextCell "C34" 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextCell "D34" '6.62'
Set-TextCell "E34" '  -7.24%  '

Set-TextCell "B35" 'PancakeSwap'
Set-TextCell "C35" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell "D35" '1.80'
Set-TextCell "E35" '  -7.42%  '

Set-TextCell "B36" 'Fetch.AI'
Set-TextCell "C36" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell "D36" '1.24'
Set-TextCell "E36" '  -8.75%  '

Set-TextCell "B37" 'EthereumClassic'
Set-TextCell "C37" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell "D37" '22.23'
Set-TextCell "E37" '  -2.14%  '

Set-TextCell "B38" 'Kaspa'
Set-TextCell "C38" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell "D38" '0.132'
Set-TextCell "E38" '  +0.20%  '

Set-TextCell "B39" 'WhiteBITCoin'
Set-TextCell "C39" 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextCell "D39" '22.22'
Set-TextCell "E39" '  -0.92%  '

Set-TextCell "B40" 'FirstDigitalUSD'
Set-TextCell "C40" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell "D40" '1.00'
Set-TextCell "E40" '  +0.26%  '

Set-TextCell "B41" 'USDe'
Set-TextCell "C41" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell "D41" '1.00'
Set-TextCell "E41" '  -0.04%  '

Set-TextCell "B42" 'PolygonEcosystemToken'
Set-TextCell "C42" 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextCell "D42" '0.360'
Set-TextCell "E42" '  -6.01%  '

Set-TextCell "B43" 'Monero'
Set-TextCell "C43" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell "D43" '147.42'
Set-TextCell "E43" '  +0.20%  '

Set-TextCell "B44" 'Stacks'
Set-TextCell "C44" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell "D44" '1.80'
Set-TextCell "E44" '  -8.94%  '

Set-TextCell "B45" 'Stellar'
Set-TextCell "C45" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell "D45" '0.133'
Set-TextCell "E45" '  +5.59%  '

Set-TextCell "B46" 'OKB'
Set-TextCell "C46" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell "D46" '43.37'
Set-TextCell "E46" '  -1.86%  '

Set-TextCell "B47" 'Hedera'
Set-TextCell "C47" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell "D47" '0.0672'
Set-TextCell "E47" '  +8.74%  '

Set-TextCell "B48" 'Mantle'
Set-TextCell "C48" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell "D48" '0.701'
Set-TextCell "E48" '  -10.70%  '

Set-TextCell "B49" 'ImmutableX'
Set-TextCell "C49" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell "D49" '1.18'
Set-TextCell "E49" '  -7.27%  '

Set-TextCell "B50" 'Aave'
Set-TextCell "C50" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell "D50" '155.36'
Set-TextCell "E50" '  -11.90%  '

Set-TextCell "B51" 'Filecoin'
Set-TextCell "C51" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell "D51" '3.92'
Set-TextCell "E51" '  -6.69%  '
